# Update crypto price/volume data per the Aug 15 2024 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.449.19"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.42%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.652.82"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.20%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.04"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.25%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.81%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.73"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.49%  "

# Row 10
$ws.Range("E10").Value = "  -2.73%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.341"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.94%  "

# Row 12
$ws.Range("E12").Value = "  +1.65%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.118.45"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.462.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.34%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.93"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000137"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.34%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.653.33"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.55%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "338.97"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.40"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.73%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.67%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.60"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.426"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.51%  "

# Row 25
$ws.Range("E25").Value = "  -1.82%  "

# Row 26
$ws.Range("E26").Value = "  +0.27%  "

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0801"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.13%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.16%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.69"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.00%  "

# Row 30
$ws.Range("E30").Value = "  +0.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.60"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.22%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.93"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "151.80"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.33%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.18"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.916"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.96%  "

# Row 36
$ws.Range("E36").Value = "  -5.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.872"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.85"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.14%  "

# Row 39
$ws.Range("E39").Value = "  -4.62%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.47%  "

# Row 41
$ws.Range("E41").Value = "  +0.41%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.609"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "275.44"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.60%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0973"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.71%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.49"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.19%  "

# Row 46
$ws.Range("E46").Value = "  -0.33%  "

# Row 47
$ws.Range("E47").Value = "  +1.62%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.046.32"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.69"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.17%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0229"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.55%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.44"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.78%  "

